# Updated the issue list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the Status column ("G") as "Completed" for the rows that were finished
$ws.Range("G6").Value = "Completed"
$ws.Range("G7").Value = "Completed"
$ws.Range("G9").Value = "Completed"
$ws.Range("G10").Value = "Completed"
$ws.Range("G11").Value = "Completed"

# Move the active selection to G6, matching the saved workbook state
$ws.Range("G6").Select()
